$d = $word.ActiveDocument

# --- Part 1: update the "This script is not capable..." paragraph ---
# Merge the two bookmark-split runs into one continuous run, append
# "in the code" directly (remove the _GoBack bookmark here - it will be
# recreated further down in the newly added content), and apply a
# yellow highlight to both the run and the paragraph mark.
$p1 = $d.Paragraphs(26)
$r1 = $p1.Range
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>This script is not capable, yet, of handling multiple items.  Its an MVP version with static references.  Some of the other required updates are commented in the code</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Part 2: append six new paragraphs describing the BOM / BOM
#     assignments / production policy steps after the final paragraph ---
$endPos = $d.Content.End
$ins = $d.Range($endPos, $endPos)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Run the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>create_BOM_table.sql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> script.  I have two scripts, at the moment (4/10/25) it' + [char]0x2019 + 's the bottom one separated by a --------- on line 43.</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">This will create a file ready to load directly into </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>coupa' + [char]0x2019 + 's</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> BOMs table.  Can be found @ </w:t></w:r><w:r><w:t>S:\Supply_Chain\Analytics\Inventory Allocation Maximization\MVP Data\input tables</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Run the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>create_BOM_assignments_table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> script.</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>This makes some modifications to the BOM script basically</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">I had to create a new table as a first step in this script to include the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>whs_co</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>de</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the products live at.  This will make sure product cant age at warehouses it doesn' + [char]0x2019 + 't exist at.</w:t></w:r></w:p>' +
'<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>I will have to make a production policy for all the products and I only want things to be able to be ' + [char]0x201C + 'produced' + [char]0x201D + ' aka aged at the warehouse they live at.</w:t></w:r></w:p>' +
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ins.InsertXML($xml2)

$d.Save()
